$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "sideboard" sheet: fill in the benchmark rows (C column first, in
#    the same order the data was originally entered, so new shared
#    strings land at the same table offsets; then the numeric columns).
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("sideboard")
$ws4.Range("C2").Value = "85.66:1"
$ws4.Range("C7").Value = "26.12:1"
$ws4.Range("C6").Value = "19.93:1"
$ws4.Range("C8").Value = "24.65:1"
$ws4.Range("C3").Value = "91.87:1"
$ws4.Range("C5").Value = "65.56:1"
$ws4.Range("C4").Value = "110.09:1"
$ws4.Range("B2").Value = 1.42
$ws4.Range("D2").Value = 0.9883
$ws4.Range("E2").Value = 0.5603
$ws4.Range("F2").Value = 700
$ws4.Range("G2").Value = 34.79
$ws4.Range("H2").Value = 0.62
$ws4.Range("I2").Value = 0.974281
$ws4.Range("J2").Value = 0.003123
$ws4.Range("K2").Value = 28
$ws4.Range("L2").Value = 3.7
$ws4.Range("B3").Value = 1.32
$ws4.Range("D3").Value = 0.9891
$ws4.Range("E3").Value = 0.5225
$ws4.Range("F3").Value = 11
$ws4.Range("G3").Value = 34.36
$ws4.Range("H3").Value = 0.56
$ws4.Range("I3").Value = 0.974023
$ws4.Range("J3").Value = 0.002848
$ws4.Range("K3").Value = 2.7
$ws4.Range("L3").Value = 2.5
$ws4.Range("B4").Value = 1.1
$ws4.Range("D4").Value = 0.9909
$ws4.Range("E4").Value = 0.436
$ws4.Range("F4").Value = 13
$ws4.Range("G4").Value = 34.55
$ws4.Range("H4").Value = 0.13
$ws4.Range("I4").Value = 0.974422
$ws4.Range("J4").Value = 0.000823
$ws4.Range("K4").Value = 11.1
$ws4.Range("L4").Value = 2.4
$ws4.Range("B5").Value = 1.85
$ws4.Range("D5").Value = 0.9847
$ws4.Range("E5").Value = 0.7321
$ws4.Range("F5").Value = 17
$ws4.Range("G5").Value = 34.53
$ws4.Range("H5").Value = 0.54
$ws4.Range("I5").Value = 0.974113
$ws4.Range("J5").Value = 0.002742
$ws4.Range("K5").Value = 7.1
$ws4.Range("L5").Value = 2.5
$ws4.Range("B6").Value = 6.1
$ws4.Range("D6").Value = 0.9498
$ws4.Range("E6").Value = 2.4081
$ws4.Range("F6").Value = 13
$ws4.Range("G6").Value = 34.14
$ws4.Range("H6").Value = 1.02
$ws4.Range("I6").Value = 0.97418
$ws4.Range("J6").Value = 0.005267
$ws4.Range("K6").Value = 3.9
$ws4.Range("L6").Value = 1.29
$ws4.Range("B7").Value = 4.65
$ws4.Range("D7").Value = 0.9617
$ws4.Range("E7").Value = 1.8379
$ws4.Range("F7").Value = 11
$ws4.Range("G7").Value = 34.62
$ws4.Range("H7").Value = 1
$ws4.Range("I7").Value = 0.974473
$ws4.Range("J7").Value = 0.004978
$ws4.Range("K7").Value = 17.32
$ws4.Range("L7").Value = 1.19
$ws4.Range("B8").Value = 4.93
$ws4.Range("D8").Value = 0.9594
$ws4.Range("E8").Value = 1.9476
$ws4.Range("F8").Value = 14
$ws4.Range("G8").Value = 34.48
$ws4.Range("H8").Value = 1
$ws4.Range("I8").Value = 0.973978
$ws4.Range("J8").Value = 0.005018
$ws4.Range("K8").Value = 8.21
$ws4.Range("L8").Value = 1.12

# ------------------------------------------------------------------
# 2) Add the new "tarot" sheet as a copy of "sideboard" placed right
#    after it, then rename it.
# ------------------------------------------------------------------
$ws4.Copy([System.Reflection.Missing]::Value, $ws4)
$ws5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5.Name = "tarot"

# ------------------------------------------------------------------
# 3) "tarot" sheet: footer labels, then benchmark rows (same
#    C-column-first ordering trick as above).
# ------------------------------------------------------------------
$ws5.Range("A19").Value = "Grid: 17 x 17"
$ws5.Range("A20").Value = "Img Size: 1024 x 1024"
$ws5.Range("A21").Value = "Orig. Size (MB): 1734.00 MB"
$ws5.Range("C2").Value = "150.56:1"
$ws5.Range("C8").Value = "31.39:1"
$ws5.Range("C6").Value = "30.30:1"
$ws5.Range("C7").Value = "37.03:1"
$ws5.Range("C4").Value = "964.52:1"
$ws5.Range("C5").Value = "366.82:1"
$ws5.Range("C3").Value = "380.96:1"
$ws5.Range("B2").Value = 11.52
$ws5.Range("D2").Value = 0.9934
$ws5.Range("E2").Value = 0.3188
$ws5.Range("F2").Value = 700
$ws5.Range("G2").Value = 36.99
$ws5.Range("H2").Value = 0.5
$ws5.Range("I2").Value = 0.97371
$ws5.Range("J2").Value = 0.002051
$ws5.Range("K2").Value = 316.5
$ws5.Range("L2").Value = 48.9
$ws5.Range("B3").Value = 4.55
$ws5.Range("D3").Value = 0.9974
$ws5.Range("E3").Value = 0.126
$ws5.Range("F3").Value = 17
$ws5.Range("G3").Value = 36.32
$ws5.Range("H3").Value = 0.84
$ws5.Range("I3").Value = 0.974317
$ws5.Range("J3").Value = 0.0022
$ws5.Range("K3").Value = 16.9
$ws5.Range("L3").Value = 19.2
$ws5.Range("B4").Value = 1.8
$ws5.Range("D4").Value = 0.999
$ws5.Range("E4").Value = 0.0498
$ws5.Range("F4").Value = 31
$ws5.Range("G4").Value = 35.97
$ws5.Range("H4").Value = 0.35
$ws5.Range("I4").Value = 0.973676
$ws5.Range("J4").Value = 0.001049
$ws5.Range("K4").Value = 54.4
$ws5.Range("L4").Value = 19.1
$ws5.Range("B5").Value = 4.73
$ws5.Range("D5").Value = 0.9973
$ws5.Range("E5").Value = 0.1309
$ws5.Range("F5").Value = 32
$ws5.Range("G5").Value = 36.19
$ws5.Range("H5").Value = 0.85
$ws5.Range("I5").Value = 0.97334
$ws5.Range("J5").Value = 0.002471
$ws5.Range("K5").Value = 35.3
$ws5.Range("L5").Value = 19.2
$ws5.Range("B6").Value = 57.23
$ws5.Range("D6").Value = 0.967
$ws5.Range("E6").Value = 1.5841
$ws5.Range("F6").Value = 13
$ws5.Range("G6").Value = 36.57
$ws5.Range("H6").Value = 2.32
$ws5.Range("I6").Value = 0.97402
$ws5.Range("J6").Value = 0.007871
$ws5.Range("K6").Value = 22.95
$ws5.Range("L6").Value = 9.58
$ws5.Range("B7").Value = 46.82
$ws5.Range("D7").Value = 0.973
$ws5.Range("E7").Value = 1.2962
$ws5.Range("F7").Value = 12
$ws5.Range("G7").Value = 36.86
$ws5.Range("H7").Value = 2.36
$ws5.Range("I7").Value = 0.973792
$ws5.Range("J7").Value = 0.007444
$ws5.Range("K7").Value = 130
$ws5.Range("L7").Value = 8.52
$ws5.Range("B8").Value = 55.23
$ws5.Range("D8").Value = 0.9681
$ws5.Range("E8").Value = 1.529
$ws5.Range("F8").Value = 15
$ws5.Range("G8").Value = 36.96
$ws5.Range("H8").Value = 2.41
$ws5.Range("I8").Value = 0.973736
$ws5.Range("J8").Value = 0.007627
$ws5.Range("K8").Value = 65
$ws5.Range("L8").Value = 9.56

# ------------------------------------------------------------------
# 4) Selections: "sideboard" -> L2 (no longer the active tab),
#    "tarot" -> G13 (becomes the active tab).
# ------------------------------------------------------------------
$ws4.Activate()
$ws4.Range("L2").Select()
$ws5.Activate()
$ws5.Range("G13").Select()
